$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old "statut" emoji markers to the new ones (column A).
$emojiMap = @{
    "🟥" = "📕"
    "⬛" = "📘"
    "🟧" = "📙"
    "🟩" = "📗"
}

# Mapping of old "statut_label" color words to the new ones (column B).
$labelMap = @{
    "noir" = "bleu"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = $cellA.Value2
    if ($emojiMap.ContainsKey($valA)) {
        $cellA.Value2 = $emojiMap[$valA]
    }

    $cellB = $ws.Cells.Item($r, 2)
    $valB = $cellB.Value2
    if ($labelMap.ContainsKey($valB)) {
        $cellB.Value2 = $labelMap[$valB]
    }
}
